$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Helper: replace the text at a specific character range with new
# text, then nudge a (no-op) formatting property on it so the engine
# keeps it as its own run instead of silently re-merging it into an
# adjacent run that happens to carry identical formatting.
# -----------------------------------------------------------------
function Set-RunRange($doc, $start, $length, $replacement) {
    $r = $doc.Range($start, $start + $length)
    $r.Text = $replacement
    $r2 = $doc.Range($start, $start + $replacement.Length)
    $r2.Bold = 1
    $r2.Bold = 0
    return $r2
}

# -----------------------------------------------------------------
# Change 1: "19 - Upgrade Badge - 86 x 86" becomes four separate runs:
#   "19 - Upgrade Badge - " / "132" / " x " / "133"
# (the badge size changes from 86 x 86 to 132 x 133, and the numbers
# end up as their own runs, same as if someone had selected each
# number in turn and typed the replacement).
# -----------------------------------------------------------------

$anchor1 = "19 - Upgrade Badge - 86 x 86"
$full = $d.Content.Text
$idx = $full.IndexOf($anchor1)
if ($idx -lt 0) {
    throw "anchor text not found: $anchor1"
}
$prefixLen = ("19 - Upgrade Badge - ").Length

# "86" -> "132"
$numStart = $idx + $prefixLen
Set-RunRange $d $numStart 2 "132" | Out-Null

# " x " stays " x " but becomes its own run
$full = $d.Content.Text
$idx = $full.IndexOf("19 - Upgrade Badge - 132 x 86")
$sepStart = $idx + ("19 - Upgrade Badge - 132").Length
Set-RunRange $d $sepStart 3 " x " | Out-Null

# "86" -> "133"
$full = $d.Content.Text
$idx = $full.IndexOf("19 - Upgrade Badge - 132 x 86")
$numStart2 = $idx + ("19 - Upgrade Badge - 132 x ").Length
Set-RunRange $d $numStart2 2 "133" | Out-Null

# -----------------------------------------------------------------
# Change 2: the runs spelling "2" + "1" + " - Upgrade " + "Close Button"
# + " - " + "39 x 36" collapse into a single run with the same text,
# "21 - Upgrade Close Button - 39 x 36" (formatting unchanged).
# A plain find & replace-in-place merges these adjacent, identically
# formatted runs into one.
# -----------------------------------------------------------------

$oldText2 = "21 - Upgrade Close Button - 39 x 36"
$d.Content.Find.Execute($oldText2, $true, $false, $false, $false, $false,
                         $true, 1, $false, $oldText2, 2) | Out-Null
